$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1.48
$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = 7.6
$ws.Range("N2").Value = 5.3
$ws.Range("T2").Value = 1.84
$ws.Range("AB2").Value = 10.5
$ws.Range("AG2").Value = 10.5
$ws.Range("AK2").Value = 15
$ws.Range("AL2").Value = 32
$ws.Range("AO2").Value = 140
$ws.Range("G3").Value = 2.04
$ws.Range("H3").Value = 3.8
$ws.Range("N3").Value = 4.9
$ws.Range("P3").Value = 2.38
$ws.Range("R3").Value = 1.54
$ws.Range("AF3").Value = 14.5
$ws.Range("AL3").Value = 46
$ws.Range("J4").Value = 3.45
$ws.Range("P4").Value = 2.16
$ws.Range("Q4").Value = 1.5
$ws.Range("F5").Value = 1.49
$ws.Range("G5").Value = 1.63
$ws.Range("H5").Value = 4.6
$ws.Range("J5").Value = 4.3
$ws.Range("G6").Value = 2.46
$ws.Range("I6").Value = 3.4
$ws.Range("N6").Value = 3.55
$ws.Range("P6").Value = 1.81
$ws.Range("Q6").Value = 2.2
$ws.Range("R6").Value = 1.31
$ws.Range("AL6").Value = 55
$ws.Range("AN6").Value = 23
$ws.Range("G7").Value = 2.6
$ws.Range("H7").Value = 2.66
$ws.Range("J7").Value = 4.1
$ws.Range("X7").Value = 28
$ws.Range("AH7").Value = 13.5
$ws.Range("AL7").Value = 28
$ws.Range("H8").Value = 1.38
$ws.Range("I8").Value = 1.39
$ws.Range("T8").Value = 1.96
$ws.Range("H9").Value = 19.5
$ws.Range("Q9").Value = 1.5
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 2.28
$ws.Range("X9").Value = 32
$ws.Range("AC9").Value = 18.5
$ws.Range("AJ9").Value = 8.6
$ws.Range("AK9").Value = 14.5
$ws.Range("AL9").Value = 140
$ws.Range("AN9").Value = 3.45
$ws.Range("G10").Value = 1.31
$ws.Range("H10").Value = 10.5
$ws.Range("I10").Value = 12
$ws.Range("P10").Value = 3.4
$ws.Range("R10").Value = 2
$ws.Range("S10").Value = 1.96
$ws.Range("T10").Value = 1.72
$ws.Range("U10").Value = 2.3
$ws.Range("Z10").Value = 140
$ws.Range("AA10").Value = 390
$ws.Range("AB10").Value = 15
$ws.Range("AC10").Value = 16.5
$ws.Range("AE10").Value = 160
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 25
$ws.Range("AJ10").Value = 12.5
$ws.Range("AL10").Value = 27
$ws.Range("AN10").Value = 3.55
$ws.Range("AO10").Value = 130
$ws.Range("K11").Value = 4.6
$ws.Range("X11").Value = 18
$ws.Range("H12").Value = 2.44
$ws.Range("I12").Value = 2.48
$ws.Range("F13").Value = 2.22
$ws.Range("G13").Value = 2.76
$ws.Range("I13").Value = 4.2
$ws.Range("J13").Value = 2.92
$ws.Range("P13").Value = 1.82
$ws.Range("Q13").Value = 1.98
